$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction0")

$values = @(1,3,4,5,6,8,9,10,11,12,14,15,16,17,19,20,21,22,23,24,25,27,29,31,32,34,35,36,37,38,39,41,42,46,47,48,49,50,51,52,53,54,55,58,59,60,61,63,64,65,69,70,72,73,74,75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
